$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Target cluster column changes from "ECs" to "MuSCs" for rows 2-4
$ws.Range("D2").Value = "MuSCs"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("D4").Value = "MuSCs"

# Row 2 updated TPM-derived values
$ws.Range("G2").Value = 7.010007333333333
$ws.Range("H2").Value = 21.030022
$ws.Range("I2").Value = 0.1776683022271353
$ws.Range("J2").Value = 0.1776683022271353
$ws.Range("M2").Value = 0.01259466666666667
$ws.Range("N2").Value = 0.037784
$ws.Range("Q2").Value = 0.08828870569422222
$ws.Range("R2").Value = 0.7945983512479999
$ws.Range("S2").Value = 0.1776683022271353
$ws.Range("T2").Value = 0.1776683022271353

# Row 3 updated TPM-derived values
$ws.Range("I3").Value = 0.6975930844911837
$ws.Range("J3").Value = 0.6975930844911837
$ws.Range("M3").Value = 0.01259466666666667
$ws.Range("N3").Value = 0.037784
$ws.Range("Q3").Value = 0.3466549168248889
$ws.Range("R3").Value = 3.119894251424
$ws.Range("S3").Value = 0.6975930844911837
$ws.Range("T3").Value = 0.6975930844911837

# Row 4 updated TPM-derived values
$ws.Range("G4").Value = 4.921635333333334
$ws.Range("H4").Value = 14.764906
$ws.Range("I4").Value = 0.124738613281681
$ws.Range("J4").Value = 0.124738613281681
$ws.Range("M4").Value = 0.01259466666666667
$ws.Range("N4").Value = 0.037784
$ws.Range("Q4").Value = 0.06198635647822222
$ws.Range("R4").Value = 0.557877208304
$ws.Range("S4").Value = 0.124738613281681
$ws.Range("T4").Value = 0.124738613281681
